$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.034373783061148
$ws.Cells.Item(2, 4).Value = 1.035181237456496
$ws.Cells.Item(2, 5).Value = 1.042369620479215
$ws.Cells.Item(2, 6).Value = 1.050605084691096
$ws.Cells.Item(2, 9).Value = 1.027244390209967
$ws.Cells.Item(2, 10).Value = 1.039493379901441
$ws.Cells.Item(2, 11).Value = 1.037978578778893
$ws.Cells.Item(2, 12).Value = 1.045146501476882
$ws.Cells.Item(2, 13).Value = 1.053358898238454
$ws.Cells.Item(2, 14).Value = 1.016869558858878
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.036094151700241
$ws.Cells.Item(3, 4).Value = 1.036721729538662
$ws.Cells.Item(3, 5).Value = 1.043891054702674
$ws.Cells.Item(3, 6).Value = 1.052196396358004
$ws.Cells.Item(3, 9).Value = 1.02724773974027
$ws.Cells.Item(3, 10).Value = 1.040852956551111
$ws.Cells.Item(3, 11).Value = 1.039326469015922
$ws.Cells.Item(3, 12).Value = 1.046476892297568
$ws.Cells.Item(3, 13).Value = 1.054760684507516
$ws.Cells.Item(3, 14).Value = 1.017342324362704
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.037205171184127
$ws.Cells.Item(4, 4).Value = 1.037716714171726
$ws.Cells.Item(4, 5).Value = 1.044873497949155
$ws.Cells.Item(4, 6).Value = 1.053223404702984
$ws.Cells.Item(4, 9).Value = 1.02724730102661
$ws.Cells.Item(4, 10).Value = 1.041730262863826
$ws.Cells.Item(4, 11).Value = 1.040196345213594
$ws.Cells.Item(4, 12).Value = 1.047335227672855
$ws.Cells.Item(4, 13).Value = 1.055664579527869
$ws.Cells.Item(4, 14).Value = 1.017646837374818
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.037671736232499
$ws.Cells.Item(5, 4).Value = 1.038134581350494
$ws.Cells.Item(5, 5).Value = 1.045286042140667
$ws.Cells.Item(5, 6).Value = 1.053654528519193
$ws.Cells.Item(5, 9).Value = 1.0272464924073
$ws.Cells.Item(5, 10).Value = 1.042098510409098
$ws.Cells.Item(5, 11).Value = 1.040561500496832
$ws.Cells.Item(5, 12).Value = 1.047695478762204
$ws.Cells.Item(5, 13).Value = 1.056043831175544
$ws.Cells.Item(5, 14).Value = 1.017774523296078
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.037750045168393
$ws.Cells.Item(6, 4).Value = 1.038204718509363
$ws.Cells.Item(6, 5).Value = 1.045355282568656
$ws.Cells.Item(6, 6).Value = 1.053726879430596
$ws.Cells.Item(6, 9).Value = 1.027246320025495
$ws.Cells.Item(6, 10).Value = 1.042160307533878
$ws.Cells.Item(6, 11).Value = 1.040622780250626
$ws.Cells.Item(6, 12).Value = 1.047755931995285
$ws.Cells.Item(6, 13).Value = 1.056107465730973
$ws.Cells.Item(6, 14).Value = 1.017795942989108
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.037211407425644
$ws.Cells.Item(7, 4).Value = 1.037722299386268
$ws.Cells.Item(7, 5).Value = 1.04487901223697
$ws.Cells.Item(7, 6).Value = 1.053229167863595
$ws.Cells.Item(7, 9).Value = 1.027247292674673
$ws.Cells.Item(7, 10).Value = 1.041735185639581
$ws.Cells.Item(7, 11).Value = 1.040201226547434
$ws.Cells.Item(7, 12).Value = 1.047340043683104
$ws.Cells.Item(7, 13).Value = 1.055669650024498
$ws.Cells.Item(7, 14).Value = 1.017648544817722
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.03495564766531
$ws.Cells.Item(8, 4).Value = 1.035702236528639
$ws.Cells.Item(8, 5).Value = 1.042884222038555
$ws.Cells.Item(8, 6).Value = 1.051143436388263
$ws.Cells.Item(8, 9).Value = 1.027246061642842
$ws.Cells.Item(8, 10).Value = 1.039953364460021
$ws.Cells.Item(8, 11).Value = 1.038434586460362
$ws.Cells.Item(8, 12).Value = 1.04559663997606
$ws.Cells.Item(8, 13).Value = 1.053833297275697
$ws.Cells.Item(8, 14).Value = 1.017029623614767
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.030963477652833
$ws.Cells.Item(9, 4).Value = 1.032128237097186
$ws.Cells.Item(9, 5).Value = 1.039353167790938
$ws.Cells.Item(9, 6).Value = 1.047447145097052
$ws.Cells.Item(9, 9).Value = 1.027223944145199
$ws.Cells.Item(9, 10).Value = 1.036794506212364
$ws.Cells.Item(9, 11).Value = 1.035303497008871
$ws.Cells.Item(9, 12).Value = 1.042504845764203
$ws.Cells.Item(9, 13).Value = 1.050572827488023
$ws.Cells.Item(9, 14).Value = 1.015928149958663
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.028289626772704
$ws.Cells.Item(10, 4).Value = 1.029735225813684
$ws.Cells.Item(10, 5).Value = 1.036987742621495
$ws.Cells.Item(10, 6).Value = 1.044968198906228
$ws.Cells.Item(10, 9).Value = 1.027195799766004
$ws.Cells.Item(10, 10).Value = 1.03467515841463
$ws.Cells.Item(10, 11).Value = 1.033203358036858
$ws.Cells.Item(10, 12).Value = 1.04042981581603
$ws.Cells.Item(10, 13).Value = 1.048382058226413
$ws.Cells.Item(10, 14).Value = 1.015186325075979
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.027128694836838
$ws.Cells.Item(11, 4).Value = 1.028696421841951
$ws.Cells.Item(11, 5).Value = 1.035960641344245
$ws.Cells.Item(11, 6).Value = 1.043891141577845
$ws.Cells.Item(11, 9).Value = 1.027180436690687
$ws.Cells.Item(11, 10).Value = 1.033754128740645
$ws.Cells.Item(11, 11).Value = 1.032290813841847
$ws.Cells.Item(11, 12).Value = 1.039527890814807
$ws.Cells.Item(11, 13).Value = 1.047429235775854
$ws.Cells.Item(11, 14).Value = 1.014863278283275
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.026696985907149
$ws.Cells.Item(12, 4).Value = 1.028310158014451
$ws.Cells.Item(12, 5).Value = 1.035578688456946
$ws.Cells.Item(12, 6).Value = 1.043490512241061
$ws.Cells.Item(12, 9).Value = 1.027174253002702
$ws.Cells.Item(12, 10).Value = 1.033411503611098
$ws.Cells.Item(12, 11).Value = 1.031951366061661
$ws.Cells.Item(12, 12).Value = 1.03919234983258
$ws.Cells.Item(12, 13).Value = 1.047074671675544
$ws.Cells.Item(12, 14).Value = 1.01474300498809
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.02678961123724
$ws.Cells.Item(13, 4).Value = 1.028393031489414
$ws.Cells.Item(13, 5).Value = 1.035660638840895
$ws.Cells.Item(13, 6).Value = 1.043576474245993
$ws.Cells.Item(13, 9).Value = 1.027175601016492
$ws.Cells.Item(13, 10).Value = 1.033485021307545
$ws.Cells.Item(13, 11).Value = 1.032024201049754
$ws.Cells.Item(13, 12).Value = 1.039264348496354
$ws.Cells.Item(13, 13).Value = 1.047150756183597
$ws.Cells.Item(13, 14).Value = 1.014768816725513
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.027093019671335
$ws.Cells.Item(14, 4).Value = 1.028664501507029
$ws.Cells.Item(14, 5).Value = 1.035929078084054
$ws.Cells.Item(14, 6).Value = 1.043858036977144
$ws.Cells.Item(14, 9).Value = 1.02717993527874
$ws.Cells.Item(14, 10).Value = 1.033725817790865
$ws.Cells.Item(14, 11).Value = 1.032262765011915
$ws.Cells.Item(14, 12).Value = 1.039500165687646
$ws.Cells.Item(14, 13).Value = 1.0473999405903
$ws.Cells.Item(14, 14).Value = 1.014853342177585
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.027279894596687
$ws.Cells.Item(15, 4).Value = 1.028831708983127
$ws.Cells.Item(15, 5).Value = 1.036094413447836
$ws.Cells.Item(15, 6).Value = 1.044031442173559
$ws.Cells.Item(15, 9).Value = 1.027182542534014
$ws.Cells.Item(15, 10).Value = 1.033874112014386
$ws.Cells.Item(15, 11).Value = 1.032409687086861
$ws.Cells.Item(15, 12).Value = 1.039645390410084
$ws.Cells.Item(15, 13).Value = 1.047553385716549
$ws.Cells.Item(15, 14).Value = 1.014905383959404
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.028366606563533
$ws.Cells.Item(16, 4).Value = 1.029804111624471
$ws.Cells.Item(16, 5).Value = 1.03705584669584
$ws.Cells.Item(16, 6).Value = 1.045039601508059
$ws.Cells.Item(16, 9).Value = 1.027196752445519
$ws.Cells.Item(16, 10).Value = 1.034736212753594
$ws.Cells.Item(16, 11).Value = 1.033263852780625
$ws.Cells.Item(16, 12).Value = 1.040489600498279
$ws.Cells.Item(16, 13).Value = 1.048445204328222
$ws.Cells.Item(16, 14).Value = 1.01520772568337
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.029047421312242
$ws.Cells.Item(17, 4).Value = 1.03041336504775
$ws.Cells.Item(17, 5).Value = 1.037658154752911
$ws.Cells.Item(17, 6).Value = 1.045671005059741
$ws.Cells.Item(17, 9).Value = 1.027204815589163
$ws.Cells.Item(17, 10).Value = 1.035276083504908
$ws.Cells.Item(17, 11).Value = 1.033798791212378
$ws.Cells.Item(17, 12).Value = 1.041018226536117
$ws.Cells.Item(17, 13).Value = 1.049003484607651
$ws.Cells.Item(17, 14).Value = 1.015396883349261
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.029444227453648
$ws.Cells.Item(18, 4).Value = 1.030768480906522
$ws.Cells.Item(18, 5).Value = 1.03800919609154
$ws.Cells.Item(18, 6).Value = 1.0460389400699
$ws.Cells.Item(18, 9).Value = 1.027209212208461
$ws.Cells.Item(18, 10).Value = 1.035590659733113
$ws.Cells.Item(18, 11).Value = 1.03411050668889
$ws.Cells.Item(18, 12).Value = 1.041326235645381
$ws.Cells.Item(18, 13).Value = 1.049328715133524
$ws.Cells.Item(18, 14).Value = 1.015507039362894
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.029579477416431
$ws.Cells.Item(19, 4).Value = 1.030889523994547
$ws.Cells.Item(19, 5).Value = 1.038128845890255
$ws.Cells.Item(19, 6).Value = 1.046164337080484
$ws.Cells.Item(19, 9).Value = 1.027210659356121
$ws.Cells.Item(19, 10).Value = 1.035697867975113
$ws.Cells.Item(19, 11).Value = 1.034216742226999
$ws.Cells.Item(19, 12).Value = 1.041431203248162
$ws.Cells.Item(19, 13).Value = 1.049439542008394
$ws.Cells.Item(19, 14).Value = 1.015544569911778
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.028974407641435
$ws.Cells.Item(20, 4).Value = 1.030348024014211
$ws.Cells.Item(20, 5).Value = 1.037593561309225
$ws.Cells.Item(20, 6).Value = 1.045603297935776
$ws.Cells.Item(20, 9).Value = 1.027203982189993
$ws.Cells.Item(20, 10).Value = 1.035218193798734
$ws.Cells.Item(20, 11).Value = 1.033741429021437
$ws.Cells.Item(20, 12).Value = 1.040961544131924
$ws.Cells.Item(20, 13).Value = 1.0489436283839
$ws.Cells.Item(20, 14).Value = 1.015376606803091
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.027003687024635
$ws.Cells.Item(21, 4).Value = 1.028584571659745
$ws.Cells.Item(21, 5).Value = 1.035850041713501
$ws.Cells.Item(21, 6).Value = 1.043775139396246
$ws.Cells.Item(21, 9).Value = 1.027178672118087
$ws.Cells.Item(21, 10).Value = 1.033654923507397
$ws.Cells.Item(21, 11).Value = 1.032192527445434
$ws.Cells.Item(21, 12).Value = 1.039430738012697
$ws.Cells.Item(21, 13).Value = 1.047326579849229
$ws.Cells.Item(21, 14).Value = 1.014828459283107
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.025761790472235
$ws.Cells.Item(22, 4).Value = 1.027473464910828
$ws.Cells.Item(22, 5).Value = 1.0347512596363
$ws.Cells.Item(22, 6).Value = 1.042622443871916
$ws.Cells.Item(22, 9).Value = 1.027159997994741
$ws.Cells.Item(22, 10).Value = 1.032669055186728
$ws.Cells.Item(22, 11).Value = 1.03121584075549
$ws.Cells.Item(22, 12).Value = 1.038465211057193
$ws.Cells.Item(22, 13).Value = 1.046306148753718
$ws.Cells.Item(22, 14).Value = 1.014482198963883
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.02642041656238
$ws.Cells.Item(23, 4).Value = 1.028062711157996
$ws.Cells.Item(23, 5).Value = 1.035333991896582
$ws.Cells.Item(23, 6).Value = 1.043233822879835
$ws.Cells.Item(23, 9).Value = 1.02717015916808
$ws.Cells.Item(23, 10).Value = 1.033191969043217
$ws.Cells.Item(23, 11).Value = 1.031733873200157
$ws.Cells.Item(23, 12).Value = 1.0389773480094
$ws.Cells.Item(23, 13).Value = 1.046847455937614
$ws.Cells.Item(23, 14).Value = 1.014665912892522
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.029007400308185
$ws.Cells.Item(24, 4).Value = 1.030377549591341
$ws.Cells.Item(24, 5).Value = 1.037622749155032
$ws.Cells.Item(24, 6).Value = 1.045633892958386
$ws.Cells.Item(24, 9).Value = 1.027204359714317
$ws.Cells.Item(24, 10).Value = 1.035244352657713
$ws.Cells.Item(24, 11).Value = 1.033767349468977
$ws.Cells.Item(24, 12).Value = 1.040987157488687
$ws.Cells.Item(24, 13).Value = 1.04897067608615
$ws.Cells.Item(24, 14).Value = 1.015385769446757
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.03199767518978
$ws.Cells.Item(25, 4).Value = 1.033053975494417
$ws.Cells.Item(25, 5).Value = 1.040267989603467
$ws.Cells.Item(25, 6).Value = 1.048405276625597
$ws.Cells.Item(25, 9).Value = 1.027232024627227
$ws.Cells.Item(25, 10).Value = 1.0376134697751
$ws.Cells.Item(25, 11).Value = 1.036115159833956
$ws.Cells.Item(25, 12).Value = 1.043306543153712
$ws.Cells.Item(25, 13).Value = 1.051418711013375
$ws.Cells.Item(25, 14).Value = 1.016214215250059
